$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the numeric threshold values in B2:C4
$ws.Range("B2").Value = 5.1
$ws.Range("C2").Value = 9
$ws.Range("B3").Value = 3.8
$ws.Range("C3").Value = 7
$ws.Range("C4").Value = 1.5

# Widen columns A and C to fit the longer parameter names, and let column B
# fall back to a plain (non-bestFit) width.
$ws.Columns.Item(1).ColumnWidth = 26.285714285714285
$ws.Columns.Item(2).ColumnWidth = 8.43
$ws.Columns.Item(3).ColumnWidth = 26.57142857142857
